# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.759.96"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "1.613.97"
$ws.Range("E3").Value = "  -1.64%  "

$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").Value = "'208.54"
$ws.Range("E5").Value = "  -2.06%  "

$ws.Range("D6").Value = "'0.519"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("D8").Value = "'23.17"
$ws.Range("E8").Value = "  -1.62%  "

$ws.Range("D9").Value = "'0.255"
$ws.Range("E9").Value = "  -2.14%  "

$ws.Range("D10").Value = "'0.0607"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("D11").Value = "'0.0873"
$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("D12").Value = "1.837.09"
$ws.Range("E12").Value = "  -1.94%  "

$ws.Range("D13").Value = "1.609.72"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").Value = "'0.557"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").Value = "'64.79"
$ws.Range("E16").Value = "  -1.69%  "

$ws.Range("D17").Value = "27.687.57"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("D18").Value = "'228.01"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").Value = "0.0₃0716"
$ws.Range("E19").Value = "  -1.02%  "

$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("D21").Value = "'0.994"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  -6.06%  "

$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  -2.71%  "

$ws.Range("D25").Value = "'153.93"
$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("D26").Value = "'6.85"
$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").Value = "'15.41"
$ws.Range("E28").Value = "  -1.80%  "

$ws.Range("D29").Value = "'0.994"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -2.34%  "

$ws.Range("D31").Value = "'0.0477"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").Value = "'3.38"
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("D33").Value = "'3.06"
$ws.Range("E33").Value = "  -2.08%  "

$ws.Range("D34").Value = "1.382.26"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  -1.21%  "

$ws.Range("D36").Value = "'0.987"
$ws.Range("E36").Value = "  +9.06%  "

$ws.Range("D37").Value = "'2.32"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").Value = "'0.0169"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'0.552"
$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("D40").Value = "'0.847"
$ws.Range("E40").Value = "  -4.01%  "

$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("D42").Value = "'0.994"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("E43").Value = "  -3.24%  "

$ws.Range("D44").Value = "'65.17"
$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("D45").Value = "'5.41"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "'2.22"
$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.747.35"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("D48").Value = "'87.37"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("D51").Value = "0.0₇0964"
$ws.Range("E51").Value = "  -8.59%  "
